$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.783.02"
$ws.Range("E2").Value = "  -1.21%  "

# Row 3
$ws.Range("D3").Value = "2.541.03"
$ws.Range("E3").Value = "  -1.82%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'310.11"
$ws.Range("E5").Value = "  -1.98%  "

# Row 6
$ws.Range("D6").Value = "'101.35"
$ws.Range("E6").Value = "  +3.78%  "

# Row 7
$ws.Range("D7").Value = "'0.572"
$ws.Range("E7").Value = "  -1.10%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  -2.23%  "

# Row 10
$ws.Range("D10").Value = "'36.37"
$ws.Range("E10").Value = "  +1.29%  "

# Row 11
$ws.Range("D11").Value = "'0.0805"
$ws.Range("E11").Value = "  -1.28%  "

# Row 12
$ws.Range("D12").Value = "'7.36"
$ws.Range("E12").Value = "  -2.44%  "

# Row 13
$ws.Range("E13").Value = "  +0.28%  "

# Row 14
$ws.Range("D14").Value = "2.927.76"
$ws.Range("E14").Value = "  -1.93%  "

# Row 15
$ws.Range("D15").Value = "'15.80"

# Row 16
$ws.Range("D16").Value = "2.506.08"
$ws.Range("E16").Value = "  -3.70%  "

# Row 17
$ws.Range("D17").Value = "'0.814"
$ws.Range("E17").Value = "  -4.07%  "

# Row 18
$ws.Range("D18").Value = "42.742.33"
$ws.Range("E18").Value = "  -1.42%  "

# Row 19
$ws.Range("D19").Value = "'6.77"
$ws.Range("E19").Value = "  -1.48%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  -1.20%  "

# Row 21
$ws.Range("D21").Value = "'12.29"
$ws.Range("E21").Value = "  -2.39%  "

# Row 22
$ws.Range("D22").Value = "'69.51"
$ws.Range("E22").Value = "  -0.22%  "

# Row 23
$ws.Range("D23").Value = "'244.84"
$ws.Range("E23").Value = "  -4.15%  "

# Row 24
$ws.Range("D24").Value = "'2.91"
$ws.Range("E24").Value = "  -2.63%  "

# Row 25
$ws.Range("D25").Value = "'2.06"
$ws.Range("E25").Value = "  -1.45%  "

# Row 26
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$ws.Range("D27").Value = "'26.22"
$ws.Range("E27").Value = "  -4.30%  "

# Row 28
$ws.Range("D28").Value = "'2.33"
$ws.Range("E28").Value = "  -4.98%  "

# Row 29
$ws.Range("D29").Value = "'39.52"
$ws.Range("E29").Value = "  -1.72%  "

# Row 30
$ws.Range("D30").Value = "'10.21"
$ws.Range("E30").Value = "  -1.37%  "

# Row 31
$ws.Range("D31").Value = "'5.79"
$ws.Range("E31").Value = "  -1.48%  "

# Row 32
$ws.Range("D32").Value = "'156.05"
$ws.Range("E32").Value = "  -0.81%  "

# Row 33
$ws.Range("E33").Value = "  +11.64%  "

# Row 34
$ws.Range("D34").Value = "'0.0795"
$ws.Range("E34").Value = "  -1.83%  "

# Row 35
$ws.Range("D35").Value = "'2.63"
$ws.Range("E35").Value = "  -2.54%  "

# Row 36
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  -5.32%  "

# Row 37
$ws.Range("D37").Value = "'3.20"
$ws.Range("E37").Value = "  -7.01%  "

# Row 38
$ws.Range("D38").Value = "'18.33"
$ws.Range("E38").Value = "  -1.98%  "

# Row 39
$ws.Range("E39").Value = "  -0.07%  "

# Row 40
$ws.Range("E40").Value = "  +0.25%  "

# Row 41
$ws.Range("E41").Value = "  +7.74%  "

# Row 42
$ws.Range("D42").Value = "'22.30"
$ws.Range("E42").Value = "  -1.44%  "

# Row 43
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "'3.33"
$ws.Range("E43").Value = "  +1.98%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("E45").Value = "  -1.53%  "

# Row 46
$ws.Range("D46").Value = "1.984.66"
$ws.Range("E46").Value = "  -1.60%  "

# Row 47
$ws.Range("D47").Value = "'8.93"
$ws.Range("E47").Value = "  -0.74%  "

# Row 48
$ws.Range("D48").Value = "2.782.88"
$ws.Range("E48").Value = "  -1.88%  "

# Row 49
$ws.Range("D49").Value = "'80.97"
$ws.Range("E49").Value = "  -3.17%  "

# Row 50
$ws.Range("E50").Value = "  -0.63%  "

# Row 51
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'72.82"
$ws.Range("E51").Value = "  -3.56%  "
